$wb = $excel.ActiveWorkbook

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("E4").Value = "2016-03-23 09:29:48"
$wsZhCn.Range("H4").Value = "2016-03-23 09:30:32"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("E4").Value = "2016-03-23 09:29:55"
$wsDeDe.Range("H4").Value = "2016-03-23 09:30:49"
